# Add the I0 (I) and IF (J) columns to the sheet, matching the existing
# header style and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: I1 = "I0", J1 = "IF" with the same style as the other headers (B1:H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2..56: column I then column J
$data = @(
    @(7, 8),
    @(6, 7),
    @(6, 8),
    @(6, 7),
    @(6, 6),
    @(8, 9),
    @(7, 8),
    @(5, 5),
    @(8, 8),
    @(5, 6),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(6, 7),
    @(6, 7),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(6, 7),
    @(6, 7),
    @(5, 6),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(5, 6),
    @(8, 9),
    @(6, 6),
    @(9, 9),
    @(7, 8),
    @(7, 7),
    @(5, 6),
    @(9, 10),
    @(7, 8),
    @(6, 7),
    @(7, 8),
    @(7, 7),
    @(6, 7),
    @(6, 7),
    @(5, 7),
    @(4, 6),
    @(4, 5),
    @(7, 8),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(6, 7),
    @(4, 5),
    @(6, 7),
    @(5, 7),
    @(5, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
